$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.716.69"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "2.532.67"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'544.74"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").Value = "'147.54"
$ws.Range("E6").Value = "  -4.63%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.581"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "2.529.81"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'5.51"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "'0.358"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "2.969.26"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "'24.65"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").Value = "59.798.60"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D18").Value = "2.514.63"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").Value = "'11.50"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'4.38"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").Value = "'327.94"
$ws.Range("D22").Value = "'0.994"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'5.86"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").Value = "'61.68"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "'0.450"
$ws.Range("E25").Value = "  -10.27%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").Value = "'7.85"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "'1.34"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.88"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.83"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "'0.996"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "'158.85"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'19.05"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "'4.54"
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("D38").Value = "'1.75"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "'6.07"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "'316.39"
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.80"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'36.78"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "'0.838"
$ws.Range("E43").Value = "  -7.37%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'0.606"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "'10.72"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'126.94"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'0.0535"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'0.0944"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").Value = "'18.76"
$ws.Range("E51").Value = "  -4.76%  "
